# Update the workbook "Avverkningsanmälningar":
#  1. Column C ("Förändrad") for all data rows (2..119) changes from
#     2023-09-23 (45192) to 2023-10-03 (45202).
#  2. Row 119 gains an explicit row height (15, custom height flag set).
#  3. A brand-new data row 120 is appended with a new case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the "Förändrad" (changed) date for every existing data row.
for ($r = 2; $r -le 119; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2. Row 119 now carries an explicit custom row height (matches row 118).
$ws.Rows.Item(119).RowHeight = 15

# 3. Append the new row (r=120) with the new case data.
$ws.Cells.Item(120, 1).Value = "A 46901-2023"

$ws.Cells.Item(120, 2).Value = 45201
$ws.Cells.Item(120, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(120, 3).Value = 45202
$ws.Cells.Item(120, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(120, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(120, 5).Value = "MULLSJÖ"

# Column F (Markägare) is intentionally left blank for this row.

$ws.Cells.Item(120, 7).Value = 1
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 13).Value = 0
$ws.Cells.Item(120, 14).Value = 0
$ws.Cells.Item(120, 15).Value = 0
$ws.Cells.Item(120, 16).Value = 0
$ws.Cells.Item(120, 17).Value = 0

# Column R keeps the wrap-text style used throughout the table, but stays
# empty (same as the row above it).
$ws.Cells.Item(120, 18).WrapText = $true
